$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 31
$ws.Cells.Item($row, 1).Value = "2025-04-28 23:15:56"
$ws.Cells.Item($row, 2).Value = 224
